# Update "detection field test data.xlsx" - append a new field-test record
# to the "dog" sheet (row 16): date 2025-05-31, type PRESENCE, target placed
# 10:40, searched 12:40.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")

# Carry the formatting (date / time number formats) down from row 15 into
# the new row 16 before writing values, same as Excel does when you type
# directly below existing formatted data.
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Value = 45808                     # Date: 31 May 2025
$ws.Range("B16").Value = "PRESENCE"                 # Type
$ws.Range("C16").Value = 0.44444444444444442        # Time target placed: 10:40
$ws.Range("D16").Value = 0.52777777777777779        # Time of search: 12:40

# Reflect where the user ended up after entering the row: scrolled back so
# column A is visible again, with the next empty cell (E16) selected.
$ws.Activate() | Out-Null
$ws.Range("E16").Select() | Out-Null
